# Update the NATMI TPM ligand-receptor table (Bmp7 -> Bmpr1b):
#  - "Resolving-Mac" is dropped as a possible "Target cluster" value, which
#    removes the two rows that targeted it (old rows 4 and 7) and, once that
#    shared string is no longer referenced, the "Resolving-Mac" entry itself.
#  - The remaining FAPs/MuSCs x FAPs/MuSCs combinations get refreshed
#    detection-rate / expression / specificity figures from the new TPM run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two rows whose "Target cluster" is "Resolving-Mac" (old sheet
# rows 7 and 4). Delete the lower row first so the row-4 index still refers
# to its original row once row 7 is gone.
$ws.Rows("7:7").Delete()
$ws.Rows("4:4").Delete()

# Row 2: FAPs -> FAPs
$ws.Cells.Item(2,1).Value2  = "FAPs"
$ws.Cells.Item(2,2).Value2  = "Bmp7"
$ws.Cells.Item(2,3).Value2  = "Bmpr1b"
$ws.Cells.Item(2,4).Value2  = "FAPs"
$ws.Cells.Item(2,5).Value2  = 1
$ws.Cells.Item(2,6).Value2  = 0.3333333333333333
$ws.Cells.Item(2,7).Value2  = 0.054589
$ws.Cells.Item(2,8).Value2  = 0.163767
$ws.Cells.Item(2,9).Value2  = 0.8203074518761176
$ws.Cells.Item(2,10).Value2 = 0.8725723693674974
$ws.Cells.Item(2,11).Value2 = 3
$ws.Cells.Item(2,12).Value2 = 1
$ws.Cells.Item(2,13).Value2 = 1.467406
$ws.Cells.Item(2,14).Value2 = 4.402218
$ws.Cells.Item(2,15).Value2 = 0.864087546066766
$ws.Cells.Item(2,16).Value2 = 0.9050919696083439
$ws.Cells.Item(2,17).Value2 = 0.08010422613399998
$ws.Cells.Item(2,18).Value2 = 0.7209380352059999
$ws.Cells.Item(2,19).Value2 = 0.7088174531119162
$ws.Cells.Item(2,20).Value2 = 0.7897582444166475

# Row 3: FAPs -> MuSCs
$ws.Cells.Item(3,1).Value2  = "FAPs"
$ws.Cells.Item(3,2).Value2  = "Bmp7"
$ws.Cells.Item(3,3).Value2  = "Bmpr1b"
$ws.Cells.Item(3,4).Value2  = "MuSCs"
$ws.Cells.Item(3,5).Value2  = 1
$ws.Cells.Item(3,6).Value2  = 0.3333333333333333
$ws.Cells.Item(3,7).Value2  = 0.054589
$ws.Cells.Item(3,8).Value2  = 0.163767
$ws.Cells.Item(3,9).Value2  = 0.8203074518761176
$ws.Cells.Item(3,10).Value2 = 0.8725723693674974
$ws.Cells.Item(3,11).Value2 = 2
$ws.Cells.Item(3,12).Value2 = 1
$ws.Cells.Item(3,13).Value2 = 0.2308085
$ws.Cells.Item(3,14).Value2 = 0.461617
$ws.Cells.Item(3,15).Value2 = 0.135912453933234
$ws.Cells.Item(3,16).Value2 = 0.09490803039165596
$ws.Cells.Item(3,17).Value2 = 0.0125996052065
$ws.Cells.Item(3,18).Value2 = 0.075597631239
$ws.Cells.Item(3,19).Value2 = 0.1114899987642014
$ws.Cells.Item(3,20).Value2 = 0.0828141249508497

# Row 4 (was row 5 before the deletions): MuSCs -> FAPs
$ws.Cells.Item(4,1).Value2  = "MuSCs"
$ws.Cells.Item(4,2).Value2  = "Bmp7"
$ws.Cells.Item(4,3).Value2  = "Bmpr1b"
$ws.Cells.Item(4,4).Value2  = "FAPs"
$ws.Cells.Item(4,5).Value2  = 1
$ws.Cells.Item(4,6).Value2  = 0.5
$ws.Cells.Item(4,7).Value2  = 0.011958
$ws.Cells.Item(4,8).Value2  = 0.023916
$ws.Cells.Item(4,9).Value2  = 0.1796925481238824
$ws.Cells.Item(4,10).Value2 = 0.1274276306325027
$ws.Cells.Item(4,11).Value2 = 3
$ws.Cells.Item(4,12).Value2 = 1
$ws.Cells.Item(4,13).Value2 = 1.467406
$ws.Cells.Item(4,14).Value2 = 4.402218
$ws.Cells.Item(4,15).Value2 = 0.864087546066766
$ws.Cells.Item(4,16).Value2 = 0.9050919696083439
$ws.Cells.Item(4,17).Value2 = 0.017547240948
$ws.Cells.Item(4,18).Value2 = 0.105283445688
$ws.Cells.Item(4,19).Value2 = 0.1552700929548498
$ws.Cells.Item(4,20).Value2 = 0.1153337251916964

# Row 5 (was row 6 before the deletions): MuSCs -> MuSCs
$ws.Cells.Item(5,1).Value2  = "MuSCs"
$ws.Cells.Item(5,2).Value2  = "Bmp7"
$ws.Cells.Item(5,3).Value2  = "Bmpr1b"
$ws.Cells.Item(5,4).Value2  = "MuSCs"
$ws.Cells.Item(5,5).Value2  = 1
$ws.Cells.Item(5,6).Value2  = 0.5
$ws.Cells.Item(5,7).Value2  = 0.011958
$ws.Cells.Item(5,8).Value2  = 0.023916
$ws.Cells.Item(5,9).Value2  = 0.1796925481238824
$ws.Cells.Item(5,10).Value2 = 0.1274276306325027
$ws.Cells.Item(5,11).Value2 = 2
$ws.Cells.Item(5,12).Value2 = 1
$ws.Cells.Item(5,13).Value2 = 0.2308085
$ws.Cells.Item(5,14).Value2 = 0.461617
$ws.Cells.Item(5,15).Value2 = 0.135912453933234
$ws.Cells.Item(5,16).Value2 = 0.09490803039165596
$ws.Cells.Item(5,17).Value2 = 0.002760008043
$ws.Cells.Item(5,18).Value2 = 0.011040032172
$ws.Cells.Item(5,19).Value2 = 0.0244224551690326
$ws.Cells.Item(5,20).Value2 = 0.01209390544080627
